# Auto-update: Pipeline state & metrics 2026-02-22T20:02:37Z
# Applies updated HDD threshold values to the "HDD Matrix" sheet (row 13 & 14).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HDD Matrix")

# Row 13: "HDDs Cold/Warmer than 10-yr Norm"
$ws.Range("B13").Value = 10
$ws.Range("W13").Value = -105
$ws.Range("Z13").Value = 2

# Row 14: "#Days with HDDs above 10yr-normals"
$ws.Range("B14").Value = 76
$ws.Range("Z14").Value = 77
